$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "agnihotriaman@gmail.com"
$ws.Range("B2").Value = "124ef1"
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = "['easy', 'medium']"
$ws.Range("E2").Value = "['Flipkart', 'Amazon']"
$ws.Range("F2").Value = "[]"
$ws.Range("G2").Value = "None"
